$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '25.829.09', '  -0.05%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.733.94', '  +0.04%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.001', '  +0.07%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '236.21', '  +2.17%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.002', '  +0.15%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.5107', '  -0.73%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2710', '  -2.22%  '),
    @(9, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '39.58', '  +0.52%  '),
    @(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06085', '  -0.31%  '),
    @(11, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.738.23', '  +0.22%  '),
    @(12, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07144', '  +1.75%  '),
    @(13, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '14.88', '  -2.15%  '),
    @(14, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6304', '  -1.76%  '),
    @(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.576', '  +1.22%  '),
    @(16, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '76.84', '  +0.08%  '),
    @(17, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.002', '  +0.16%  '),
    @(18, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.002', '  +0.07%  '),
    @(19, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '25.837.48', '  +0.03%  '),
    @(20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '11.66', '  +1.36%  '),
    @(21, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000006675', '  +0.76%  '),
    @(22, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '1.959.94', '  -0.32%  '),
    @(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.251', '  +2.62%  '),
    @(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.586', '  -1.38%  '),
    @(25, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '5.205', '  +1.71%  '),
    @(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '139.37', '  -0.38%  '),
    @(27, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.512', '  -0.23%  '),
    @(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '15.12', '  +0.67%  '),
    @(29, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.754', '  -2.19%  '),
    @(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '104.88', '  +2.85%  '),
    @(31, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '3.878', '  +5.04%  '),
    @(32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08329', '  +0.33%  '),
    @(33, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.565', '  +4.17%  '),
    @(34, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04533', '  +1.26%  '),
    @(35, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.630', '  +0.50%  '),
    @(36, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '0.9767', '  -0.30%  '),
    @(37, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.6177', '  +0.74%  '),
    @(38, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.684', '  +1.70%  '),
    @(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01586', '  +0.46%  '),
    @(40, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.002', '  +0.22%  '),
    @(41, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.905', '  -1.26%  '),
    @(42, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '97.11', '  -3.31%  '),
    @(43, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.3825', '  +0.31%  '),
    @(44, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.7323', '  +0.43%  '),
    @(45, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '4.911', '  -0.94%  '),
    @(46, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1124', '  +0.22%  '),
    @(47, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05270', '  -2.08%  '),
    @(48, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.160', '  -1.19%  '),
    @(49, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '54.33', '  +2.65%  '),
    @(50, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '30.34', '  +1.09%  '),
    @(51, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '7.535', '  -1.23%  '),
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 2).Style = "Normal"

    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 3).Style = "Normal"

    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 4).Style = "Normal"

    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 5).Style = "Normal"
}
